$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider (J_0_g constraints) ---
$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
$wsLider.Range("A2").Value = "2.09 - x"
$wsLider.Range("B2").Value = -3.09
$wsLider.Range("D2").Value = 0.86

$wsLider.Range("A3").Value = "-2.09 + x"
$wsLider.Range("B3").Value = 1.0899999999999999
$wsLider.Range("D3").Value = 0.62

$wsLider.Range("A4").Value = "41.02289999999999 + x - y - 9(x^2)"
$wsLider.Range("B4").Value = -40.02289999999999
$wsLider.Range("D4").Value = 0.58

# --- Restricciones_del_follower (Gamma restrictions) ---
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsFollower.Range("A2").Value = "-22.9596 + (-0.5 + x)*(y^2)"
$wsFollower.Range("B2").Value = 22.9596
$wsFollower.Range("D2").Value = 0.69
$wsFollower.Range("E2").Value = 0
$wsFollower.Range("F2").Value = 6.7

$wsFollower.Range("A3").Value = "-3.8 + y"
$wsFollower.Range("B3").Value = 2.8
$wsFollower.Range("D3").Value = 0.65
$wsFollower.Range("E3").Value = 4.5
$wsFollower.Range("F3").Value = 0

$wsFollower.Range("A4").Value = "-5.8 - y"
$wsFollower.Range("B4").Value = -4.8
$wsFollower.Range("D4").Value = 0.32
$wsFollower.Range("E4").Value = 3.2
$wsFollower.Range("F4").Value = 3.3000000000000003

# --- Punto_modificado (x, y) ---
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2").Value = 2.09
$wsPunto.Range("B2").Value = 3.8

# --- Vector_bf (sheet index 5) ---
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2").Value = -9.667959999999997

# --- Vector_BF (sheet index 6; name lookup is case-insensitive and would
#     collide with "Vector_bf" above, so address this sheet by position) ---
$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2").Value = 20.479599999999994
$wsBF.Range("A3").Value = -0.7199999999999999
